$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: cut the JSON payload (with its formatting) from B9 over to C9,
# then type the test-case number into the now-empty B9.
$ws.Range("B9").Cut($ws.Range("C9"))
$ws.Range("B9").Value = 6

# Row 10: same move - cut B10's JSON payload into C10, then set B10 to 7.
$ws.Range("B10").Cut($ws.Range("C10"))
$ws.Range("B10").Value = 7

# Leave the selection where the user last clicked before saving.
$ws.Range("C15").Select()
